$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert three new blank paragraphs right after "Episode One" and before
#    the "Concept game..." paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Concept game, what I want it to do, what I want it to feel like.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(1)   # wdCollapseStart
for ($i = 0; $i -lt 3; $i++) {
    $rng.InsertParagraphBefore()
}

# ---------------------------------------------------------------------------
# 2. The "Episode Two" heading used to carry a (rendering-only) page-break
#    marker directly on its run. Re-assigning the paragraph's text rebuilds
#    the run cleanly without that marker, leaving a plain text run behind -
#    matching the new layout where the break no longer lands on this line.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Episode Two", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$epTwoPara = $rng2.Paragraphs.Item(1)
$epTwoPara.Range.Text = "Episode Two"

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
